# Daily attendance processing - 2025-12-03 09:55:06
# Applies the session-analysis refresh: reordered "Recorded By" email lists,
# updated Group/Class statistics counters, and the two sessions (B1/MICROBIOLOGY
# session 1, B3/PHYSIOLOGY session 1) that flipped from "Not Recorded" to
# "Recorded" now that attendance has come in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: several "xx.x%" figures in this sheet are stored as literal TEXT
# (not numeric percentages). Assigning a "%"-suffixed string straight to
# .Value makes Excel auto-coerce it into a real percentage number (and mints
# a new number-format style in the process). Writing it as a formula that
# evaluates to the text, then collapsing the formula down to its value via a
# self Copy/PasteSpecial(values), keeps it as plain text on the original
# cell style.
function Set-TextPercent {
    param($addr, $text)
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# "Recorded By" (column G) email-list reorderings
# ---------------------------------------------------------------------------
$ws.Range("G2").Value  = "servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G18").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G24").Value = "servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G40").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G58").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G62").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G80").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G84").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G96").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G106").Value = "nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G118").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G120").Value = "amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G128").Value = "nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Class Statistics block (K4:L10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 24      # Recorded Sessions
$ws.Range("L7").Value = 10      # Missing Sessions
Set-TextPercent "L9"  "13.6%"   # Coverage %
Set-TextPercent "L10" "30.6%"   # Average Attendance %

# ---------------------------------------------------------------------------
# Group Statistics block (K14:S...) - group B1 (row 19) and B3 (row 21)
# ---------------------------------------------------------------------------
$ws.Range("O19").Value = 3
$ws.Range("P19").Value = 3
Set-TextPercent "R19" "13.6%"
Set-TextPercent "S19" "26.8%"

$ws.Range("O21").Value = 3
$ws.Range("P21").Value = 1
Set-TextPercent "R21" "13.6%"
Set-TextPercent "S21" "20.2%"

# ---------------------------------------------------------------------------
# Row 98: B1 / MICROBIOLOGY session 1 flips from "Not Recorded" to "Recorded".
# Copy the normal (non-highlighted) row formatting from row 2 first, then
# fill in the now-available attendance data.
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A98:I98").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G98").Value = "amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("H98").Value = "71/154"
$ws.Range("I98").Value = "Recorded"

# ---------------------------------------------------------------------------
# Row 150: B3 / PHYSIOLOGY session 1 flips from "Not Recorded" to "Recorded".
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A150:I150").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G150").Value = "naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("H150").Value = "46/224"
$ws.Range("I150").Value = "Recorded"
